# Cambios listos de Bodegaje e IVA
# Insert a new "Total IVA :" row (row 20) right before the existing
# "Saldo Inicial / Total Cierre / Total Ingresos" totals block, pushing the
# remaining rows (old 20-26) down to (21-27).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20; Excel inherits formatting from the row
# above (row 19 - "Total Ingresos :"), which is exactly the style the new
# "Total IVA :" row should use.
$ws.Rows.Item(20).Insert()

# Match the row height used by the surrounding total rows (60pt) before
# typing the label, so AutoFit doesn't shrink it to the text's natural size.
$ws.Rows.Item(20).RowHeight = 60

# Label for the new row.
$ws.Range("A20").Value = "Total IVA :"

# Same merge layout as the other total rows directly above it:
# label spans A:B, the (empty) value spans C:D.
$ws.Range("A20:B20").Merge()
$ws.Range("C20:D20").Merge()

# Leave the selection on the new value cell, matching the saved view state.
$ws.Range("C20:D20").Select()
